# KUCE_placements.xlsx — "checked and fixed the ers" edit
#
# 1. Company names in column B for rows that are re-appearances of an
#    already-used company get a distinguishing suffix (TCS-2, WIPRO-2, TCS-3)
#    so the duplicate-detector used below treats them as unique text while the
#    Sl.NO numbering (col A) is renumbered sequentially 1..8.
# 2. A "Highlight Duplicate Values" conditional format is applied to column B
#    (the whole column, Excel's default range when you run the gallery rule
#    from the Home ribbon).
# 3. The last active selection left on the sheet is C11 (cosmetic — matches
#    where the user clicked after finishing the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- renumber Sl.NO (column A) sequentially ------------------------------
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A9").Value = 8

# --- disambiguate repeated company names (column B) ----------------------
# (set in shared-string insertion order: TCS-2, WIPRO-2, TCS-3)
$ws.Range("B9").Value = "TCS-2"
$ws.Range("B7").Value = "WIPRO-2"
$ws.Range("B6").Value = "TCS-3"

# --- Highlight Duplicate Values conditional formatting on column B -------
$col = $ws.Range("B1:B1048576")
$fc = $col.FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1
$fc.Font.Color = 393372        # RGB(0x9C,0x00,0x06) = FF9C0006 (dark red)
$fc.Interior.Color = 13551615  # RGB(0xFF,0xC7,0xCE) = FFFFC7CE (light red)

# --- leave the same cell selected as in the authored workbook ------------
$ws.Range("C11").Select()
